# Weekly update: insert the newest week's 3 rows (Especial/Primera/Segunda)
# at the top of the data block (row 11), pushing all existing rows down by
# three. This matches the "Fruta / hortaliza, semanal" (weekly) refresh
# pattern used throughout this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows right before the current first data row (row 11).
$ws.Rows("11:13").Insert()

# --- Row 11: Especial, Provincia de Limarí ---
$ws.Range("A11").Value = 3
$ws.Range("B11").Value = "Femacal de La Calera"
$ws.Range("C11").Value = "Coquimbo"
$ws.Range("D11").Value = "2021-11-03"
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100107
$ws.Range("H11").Value = "Otros"
$ws.Range("I11").Value = 100107002
$ws.Range("J11").Value = "Chirimoya"
$ws.Range("K11").Value = "Cultivar IV Región"
$ws.Range("L11").Value = "Especial"
$ws.Range("M11").Value = 50
$ws.Range("N11").Value = 27000
$ws.Range("O11").Value = 27000
$ws.Range("P11").Value = 27000
$ws.Range("Q11").Value = "$/bandeja 10 kilos"
$ws.Range("R11").Value = "Provincia de Limarí"
$ws.Range("S11").Value = 2700
$ws.Range("T11").Value = 10

# --- Row 12: Primera, Provincia de Limarí ---
$ws.Range("A12").Value = 3
$ws.Range("B12").Value = "Femacal de La Calera"
$ws.Range("C12").Value = "Coquimbo"
$ws.Range("D12").Value = "2021-11-03"
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100107
$ws.Range("H12").Value = "Otros"
$ws.Range("I12").Value = 100107002
$ws.Range("J12").Value = "Chirimoya"
$ws.Range("K12").Value = "Cultivar IV Región"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 65
$ws.Range("N12").Value = 25000
$ws.Range("O12").Value = 25000
$ws.Range("P12").Value = 25000
$ws.Range("Q12").Value = "$/bandeja 10 kilos"
$ws.Range("R12").Value = "Provincia de Limarí"
$ws.Range("S12").Value = 2500
$ws.Range("T12").Value = 10

# --- Row 13: Segunda, Provincia de Limarí ---
$ws.Range("A13").Value = 3
$ws.Range("B13").Value = "Femacal de La Calera"
$ws.Range("C13").Value = "Coquimbo"
$ws.Range("D13").Value = "2021-11-03"
$ws.Range("E13").Value = 5
$ws.Range("F13").Value = "Fruta"
$ws.Range("G13").Value = 100107
$ws.Range("H13").Value = "Otros"
$ws.Range("I13").Value = 100107002
$ws.Range("J13").Value = "Chirimoya"
$ws.Range("K13").Value = "Cultivar IV Región"
$ws.Range("L13").Value = "Segunda"
$ws.Range("M13").Value = 60
$ws.Range("N13").Value = 22000
$ws.Range("O13").Value = 22000
$ws.Range("P13").Value = 22000
$ws.Range("Q13").Value = "$/bandeja 10 kilos"
$ws.Range("R13").Value = "Provincia de Limarí"
$ws.Range("S13").Value = 2200
$ws.Range("T13").Value = 10
